# Remove the first three data rows (2007, 2008, 2009) from the table.
# This shifts the remaining years (2010-2018) up so they occupy rows 2-10,
# matching the new dimension A1:M10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Delete()
